# Add a new "2022" column (P) to the table, mirroring the existing
# 2021 column (O) for formatting, then fill in the new data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: year header (2021 -> 2022)
$ws.Range("O3").Copy($ws.Range("P3"))
$ws.Range("P3").Value = 2022

# Row 4: count of reporting companies (14 -> 15)
$ws.Range("O4").Copy($ws.Range("P4"))
$ws.Range("P4").Value = 15

# Row 5: insurance premiums, mln soms (1252.8 -> 2130.4)
$ws.Range("O5").Copy($ws.Range("P5"))
$ws.Range("P5").Value = 2130.4

# Move the active selection to P6, matching the saved view state.
$ws.Range("P6").Select()
